$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the refreshed cryptocurrency market data.
# Numeric-looking text values (e.g. "1.00", "0.425") need to stay as literal
# text (matching the sheet's existing inline-string cells), so we force the
# Text number format before assignment and then clear it again so the cell
# keeps its original (unstyled) formatting.

$ws.Range("D2").Value = '67.404.65'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '3.516.55'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '610.39'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.56'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.15%  '
$ws.Range("D7").Value = '3.515.83'
$ws.Range("E7").Value = '  -1.22%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.425'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("D14").Value = '4.111.36'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.78'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = '3.510.58'
$ws.Range("E16").Value = '  -1.38%  '
$ws.Range("D17").Value = '67.406.82'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.23'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '443.48'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.25'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.625'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.24%  '
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  +10.86%  '
$ws.Range("D26").Value = '3.657.31'
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.21'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.03%  '
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("E30").Value = '  -2.50%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.54'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.163'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.05%  '
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.13'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.33%  '
$ws.Range("D36").Value = '3.511.23'
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("E37").Value = '  -3.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.02'
$ws.Range("D38").ClearFormats()
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '179.21'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.17'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0873'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.52%  '
$ws.Range("E44").Value = '  -3.30%  '
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.56'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.67'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("E48").Value = '  +4.73%  '
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.57'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.87%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.994'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.41%  '
